$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to remain plain text so Excel does not
# reinterpret values like "1.001" or "237.55" as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.235.39'
$ws.Range("E2").Value = '  +0.36%  '

$ws.Range("D3").Value = '1.859.31'
$ws.Range("E3").Value = '  +0.52%  '

$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '0.7015'
$ws.Range("E5").Value = '  -0.37%  '

$ws.Range("D6").Value = '237.55'
$ws.Range("E6").Value = '  -0.23%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = '0.08284'
$ws.Range("E8").Value = '  +10.30%  '

$ws.Range("D9").Value = '0.3034'
$ws.Range("E9").Value = '  -0.16%  '

$ws.Range("E10").Value = '  -0.60%  '

$ws.Range("D11").Value = '0.08180'
$ws.Range("E11").Value = '  +0.66%  '

$ws.Range("D12").Value = '1.869.92'
$ws.Range("E12").Value = '  +1.32%  '

$ws.Range("D13").Value = '5.172'
$ws.Range("E13").Value = '  -0.81%  '

$ws.Range("D14").Value = '0.7121'
$ws.Range("E14").Value = '  -1.89%  '

$ws.Range("D15").Value = '89.10'
$ws.Range("E15").Value = '  +0.17%  '

$ws.Range("D16").Value = '29.256.63'
$ws.Range("E16").Value = '  +0.74%  '

$ws.Range("D17").Value = '5.777'
$ws.Range("E17").Value = '  +0.06%  '

$ws.Range("D18").Value = '0.000007847'
$ws.Range("E18").Value = '  +2.47%  '

$ws.Range("D19").Value = '13.35'
$ws.Range("E19").Value = '  +2.22%  '

$ws.Range("D20").Value = '237.03'
$ws.Range("E20").Value = '  -0.59%  '

$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.09%  '

$ws.Range("D22").Value = '2.112.45'
$ws.Range("E22").Value = '  +1.67%  '

$ws.Range("E23").Value = '  +0.10%  '

$ws.Range("D24").Value = '7.440'
$ws.Range("E24").Value = '  -1.38%  '

$ws.Range("D25").Value = '161.78'
$ws.Range("E25").Value = '  -0.15%  '

$ws.Range("D26").Value = '8.969'
$ws.Range("E26").Value = '  -0.15%  '

$ws.Range("D27").Value = '0.1443'
$ws.Range("E27").Value = '  -0.99%  '

$ws.Range("D28").Value = '18.10'
$ws.Range("E28").Value = '  +0.44%  '

$ws.Range("E29").Value = '  +1.31%  '

$ws.Range("D30").Value = '1.436'
$ws.Range("E30").Value = '  +3.63%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '1.484'
$ws.Range("E31").Value = '  -0.62%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '4.400'
$ws.Range("E32").Value = '  -3.12%  '

$ws.Range("D33").Value = '4.058'
$ws.Range("E33").Value = '  +1.75%  '

$ws.Range("D34").Value = '0.05207'
$ws.Range("E34").Value = '  +1.18%  '

$ws.Range("E35").Value = '  -1.57%  '

$ws.Range("D36").Value = '0.7071'
$ws.Range("E36").Value = '  +1.02%  '

$ws.Range("D37").Value = '1.004'
$ws.Range("E37").Value = '  -3.32%  '

$ws.Range("D38").Value = '2.669'
$ws.Range("E38").Value = '  +0.95%  '

$ws.Range("D39").Value = '0.01845'
$ws.Range("E39").Value = '  -1.50%  '

$ws.Range("D40").Value = '2.724'
$ws.Range("E40").Value = '  +1.71%  '

$ws.Range("D41").Value = '0.9182'
$ws.Range("E41").Value = '  -2.72%  '

$ws.Range("D42").Value = '1.134.38'
$ws.Range("E42").Value = '  +5.00%  '

$ws.Range("D43").Value = '5.936'
$ws.Range("E43").Value = '  -0.82%  '

$ws.Range("D44").Value = '0.4278'
$ws.Range("E44").Value = '  -0.26%  '

$ws.Range("D45").Value = '70.65'
$ws.Range("E45").Value = '  +1.18%  '

$ws.Range("D46").Value = '0.9999'
$ws.Range("E46").Value = '  -0.03%  '

$ws.Range("D47").Value = '102.36'
$ws.Range("E47").Value = '  +0.04%  '

$ws.Range("D48").Value = '1.770'
$ws.Range("E48").Value = '  +1.51%  '

$ws.Range("D49").Value = '2.009.46'
$ws.Range("E49").Value = '  +1.78%  '

$ws.Range("D50").Value = '9.178'
$ws.Range("E50").Value = '  +0.34%  '

$ws.Range("D51").Value = '6.978'
$ws.Range("E51").Value = '  -0.88%  '
